$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("B2").Value = 7206
$ws.Range("C3").Value = 178117
$ws.Range("C4").Value = 168073
$ws.Range("C8").Value = 64.84
